$d = $word.ActiveDocument

# En dash character used throughout the list ("17 <en dash> Road city ...")
$dash = [char]8211

# Locate the last paragraph in the document (ends with "17 - Road city crossroad vertical")
$lastPara = $d.Paragraphs.Last

# Insert three new paragraphs after it, mirroring the formatting (en-GB language) of the
# preceding paragraph, and carrying the new tile descriptions.
$lastPara.Range.InsertParagraphAfter()
$p18 = $d.Paragraphs.Last
$p18.Range.Text = "18 $dash Grass trees"

$p18.Range.InsertParagraphAfter()
$p19 = $d.Paragraphs.Last
$p19.Range.Text = "19 $dash Grass trees vertical"

$p19.Range.InsertParagraphAfter()
$p20 = $d.Paragraphs.Last
$p20.Range.Text = "20 $dash Grass tree lines horizontal"
